$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5962  # was 5942
$ws.Range("L3").Value = 6493  # was 6472
$ws.Range("J4").Value = 1722  # was 1721
$ws.Range("L4").Value = 1593  # was 1592
$ws.Range("L5").Value = 386  # was 383
$ws.Range("L6").Value = 5342  # was 5325
$ws.Range("J7").Value = 26202  # was 26201
$ws.Range("L7").Value = 19776  # was 19714

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 174  # was 173
$ws.Range("L5").Value = 73  # was 72
$ws.Range("L7").Value = 641  # was 638
$ws.Range("L8").Value = 1306  # was 1300
$ws.Range("L9").Value = 113  # was 111
$ws.Range("L11").Value = 329  # was 326
$ws.Range("L20").Value = 503  # was 498
$ws.Range("L21").Value = 61  # was 60
$ws.Range("L22").Value = 62  # was 61
$ws.Range("L24").Value = 57  # was 56
$ws.Range("L25").Value = 118  # was 117
$ws.Range("L28").Value = 6  # was 7
$ws.Range("L29").Value = 1114  # was 1110
$ws.Range("L33").Value = 889  # was 886
$ws.Range("L36").Value = 249  # was 248
$ws.Range("L37").Value = 749  # was 748
$ws.Range("L43").Value = 148  # was 146
$ws.Range("L52").Value = 415  # was 413
$ws.Range("L53").Value = 218  # was 217
$ws.Range("L54").Value = 432  # was 425
$ws.Range("L55").Value = 205  # was 204
$ws.Range("L59").Value = 34  # was 33
$ws.Range("K63").Value = 175  # was 174
$ws.Range("L63").Value = 56  # was 58
$ws.Range("L65").Value = 387  # was 385
$ws.Range("L67").Value = 680  # was 679
$ws.Range("L76").Value = 301  # was 300
$ws.Range("L79").Value = 550  # was 548
$ws.Range("J83").Value = 526  # was 525
$ws.Range("K83").Value = 529  # was 530
$ws.Range("L83").Value = 431  # was 430
$ws.Range("L84").Value = 190  # was 189
$ws.Range("L85").Value = 980  # was 977
$ws.Range("L91").Value = 266  # was 264
$ws.Range("L93").Value = 100  # was 99
$ws.Range("L94").Value = 248  # was 246
$ws.Range("L95").Value = 280  # was 279
$ws.Range("L96").Value = 221  # was 222
$ws.Range("L100").Value = 37  # was 36
$ws.Range("J101").Value = 26202  # was 26201
$ws.Range("L101").Value = 19776  # was 19714

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 65  # was 66
$ws.Range("L7").Value = 221  # was 222

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 219  # was 216
$ws.Range("L7").Value = 641  # was 638

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 124  # was 122
$ws.Range("L6").Value = 81  # was 80
$ws.Range("L7").Value = 329  # was 326

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 293  # was 292
$ws.Range("L3").Value = 406  # was 405
$ws.Range("L6").Value = 205  # was 204
$ws.Range("L7").Value = 980  # was 977

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 130  # was 129
$ws.Range("L6").Value = 118  # was 117
$ws.Range("L7").Value = 415  # was 413

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 74  # was 73
$ws.Range("L7").Value = 218  # was 217

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 393  # was 390
$ws.Range("L3").Value = 461  # was 459
$ws.Range("L6").Value = 320  # was 319
$ws.Range("L7").Value = 1306  # was 1300

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 174  # was 173
$ws.Range("J4").Value = 24  # was 23
$ws.Range("K4").Value = 23  # was 24
$ws.Range("J7").Value = 526  # was 525
$ws.Range("K7").Value = 529  # was 530
$ws.Range("L7").Value = 431  # was 430

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 259  # was 256
$ws.Range("L7").Value = 889  # was 886

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 89  # was 88
$ws.Range("L7").Value = 280  # was 279

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 226  # was 225
$ws.Range("L7").Value = 749  # was 748

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 127  # was 126
$ws.Range("L6").Value = 93  # was 92
$ws.Range("L7").Value = 387  # was 385

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L5").Value = 19  # was 18
$ws.Range("L7").Value = 680  # was 679

# South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 64  # was 63
$ws.Range("L7").Value = 190  # was 189

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 109  # was 104
$ws.Range("L5").Value = 3  # was 2
$ws.Range("L6").Value = 209  # was 208
$ws.Range("L7").Value = 432  # was 425

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 329  # was 326
$ws.Range("L6").Value = 273  # was 272
$ws.Range("L7").Value = 1114  # was 1110

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L5").Value = 9  # was 8
$ws.Range("L7").Value = 301  # was 300

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 69  # was 68
$ws.Range("L7").Value = 205  # was 204

# Dunning
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L2").Value = 24  # was 23
$ws.Range("L7").Value = 57  # was 56

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 91  # was 89
$ws.Range("L7").Value = 266  # was 264

# Chinatown
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 32  # was 31
$ws.Range("L7").Value = 61  # was 60

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 178  # was 177
$ws.Range("L6").Value = 148  # was 147
$ws.Range("L7").Value = 550  # was 548

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 158  # was 156
$ws.Range("L3").Value = 175  # was 173
$ws.Range("L6").Value = 118  # was 117
$ws.Range("L7").Value = 503  # was 498

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 88  # was 87
$ws.Range("L7").Value = 249  # was 248

# West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 31  # was 30
$ws.Range("L7").Value = 100  # was 99

# Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 22  # was 21
$ws.Range("L7").Value = 37  # was 36

# West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 60  # was 59
$ws.Range("L3").Value = 60  # was 59
$ws.Range("L7").Value = 248  # was 246

# East Side
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 55  # was 54
$ws.Range("L7").Value = 118  # was 117

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 44  # was 43
$ws.Range("L6").Value = 27  # was 26
$ws.Range("L7").Value = 113  # was 111

# Montclare
$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L2").Value = 11  # was 10
$ws.Range("L7").Value = 34  # was 33

# Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L4").Value = 15  # was 14
$ws.Range("L7").Value = 174  # was 173

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 33  # was 32
$ws.Range("L7").Value = 73  # was 72

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 50  # was 48
$ws.Range("L7").Value = 148  # was 146

# Clearing
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L4").Value = 8  # was 7
$ws.Range("L7").Value = 62  # was 61

# Edison Park
$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("L4").ClearContents()  # was 1
$ws.Range("L7").Value = 6  # was 7
